# Append new Bathroom sensor-log rows (2026-01-28 afternoon readings)
# to the PIR, Humidity and Temperature sheets.
$wb = $excel.ActiveWorkbook

# --- PIR: append rows 27-39 ---
$ws = $wb.Worksheets.Item("PIR")

$pirData = @(
  @("2026-01-28","16:12:38","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:12:39","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:12:43","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:12:48","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:12:53","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:12:58","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:03","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:08","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:13","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:18","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:23","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:28","16:00","Bathroom","No Motion","Inactive"),
  @("2026-01-28","16:13:33","16:00","Bathroom","No Motion","Inactive")
)

$ws.Range("A27:A39").NumberFormat = "@"
$ws.Range("E27:E39").NumberFormat = "@"

$r = 27
foreach ($row in $pirData) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# --- Humidity: append rows 26-39 ---
$ws = $wb.Worksheets.Item("Humidity")

$humidityData = @(
  @("2026-01-28","16:12:38","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:12:39","16:00","Bathroom","87.4%","Active"),
  @("2026-01-28","16:12:41","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:12:45","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:12:49","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:12:53","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:12:57","16:00","Bathroom","87.3%","Active"),
  @("2026-01-28","16:13:01","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:13:05","16:00","Bathroom","88.3%","Active"),
  @("2026-01-28","16:13:14","16:00","Bathroom","88.2%","Active"),
  @("2026-01-28","16:13:18","16:00","Bathroom","87.4%","Active"),
  @("2026-01-28","16:13:22","16:00","Bathroom","88.2%","Active"),
  @("2026-01-28","16:13:30","16:00","Bathroom","87.3%","Active"),
  @("2026-01-28","16:13:34","16:00","Bathroom","88.2%","Active")
)

$ws.Range("A26:A39").NumberFormat = "@"
$ws.Range("E26:E39").NumberFormat = "@"

$r = 26
foreach ($row in $humidityData) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

# --- Temperature: append rows 26-39 ---
$ws = $wb.Worksheets.Item("Temperature")

$temperatureData = @(
  @("2026-01-28","16:12:39","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:12:39","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:12:42","16:00","Bathroom","22.7C","Active"),
  @("2026-01-28","16:12:45","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:12:50","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:12:54","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:12:58","16:00","Bathroom","22.7C","Active"),
  @("2026-01-28","16:13:02","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:13:06","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:13:14","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:13:18","16:00","Bathroom","22.8C","Active"),
  @("2026-01-28","16:13:22","16:00","Bathroom","22.7C","Active"),
  @("2026-01-28","16:13:30","16:00","Bathroom","22.7C","Active"),
  @("2026-01-28","16:13:34","16:00","Bathroom","22.7C","Active")
)

$ws.Range("A26:A39").NumberFormat = "@"
$ws.Range("E26:E39").NumberFormat = "@"

$r = 26
foreach ($row in $temperatureData) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}
